$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Issue #57: genre becomes required (PBCore controlled vocabulary). The
# spreadsheet fixture gains a new "Genre" column (U) with a header matching
# the existing "Genre" label already used elsewhere on the sheet, and a
# genre value ("Aviation") supplied for each of the two data rows so the
# fixture keeps validating once genre is mandatory.
$ws.Range("U2").Value = "Genre"
$ws.Range("U3").Value = "Aviation"
$ws.Range("U4").Value = "Aviation"

# Leave the sheet selection where the author ended up after adding the
# column (one row below the last data row, in the new column).
$ws.Range("U5").Select()
